$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Tyler Anderson's row (row 2): new "Last.Updated" date and new injury details text
$ws.Range("C2").Value = "August 21 2017"
$ws.Range("E2").Value = "Anderson has been transferred to the 60-day disabled list due to arthroscopic left knee surgery and is expected to be sidelined until the start of September."

# Add new row 4 for Ryan Hanigan (entered C:E first, then A:B, to mirror original authoring order)
$ws.Range("C4").Value = "August 22 2017"
$ws.Range("D4").Value = "Groin"
$ws.Range("E4").Value = "Hanigan is on the 10-day disabled list with a left groin strain and will miss an undetermined amount of game action."
$ws.Range("A4").Value = "Ryan Hanigan"
$ws.Range("B4").Value = "hanigry01"

# Update selection to match the final cursor position recorded in the workbook
$ws.Range("A16").Select()
